# "Add files via upload" - refresh the map-tracking data (dados_mapa.xlsx)
# with updated, accent-corrected headers and a new data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): fix missing Portuguese accents ---
$ws.Range("C1").Value = "Número endereço"
$ws.Range("D1").Value = "Atualização"

# --- Data row (row 2): new date + numeric values instead of text ---
$ws.Cells.Item(2, 1).Value = 45590          # Data -> 2024-10-25 (serial date)
$ws.Range("A2").NumberFormat = "YYYY-MM-DD"
$ws.Range("B2").Value = 1                    # Mapa Selecionado -> 1
$ws.Range("C2").Value = 3                    # Numero endereco -> 3
$ws.Range("D2").Value = "Mudou-se"           # Atualizacao -> Mudou-se

$ws.Range("A1").Select()
